# Auto-generated edit script: applies updated market-price derived values
# across the eight crafting-job profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 536.125
$ws.Range("I19").Value = 498.16666
$ws.Range("K19").Value = 498.16666
$ws.Range("M19").Value = -323.16666
$ws.Range("H40").Value = 1496.25
$ws.Range("I40").Value = 1481.4286
$ws.Range("J40").Value = 1600
$ws.Range("K40").Value = 1481.4286
$ws.Range("L40").Value = 1600
$ws.Range("M40").Value = -1306.4286
$ws.Range("N40").Value = -1950
$ws.Range("H64").Value = 3939.9
$ws.Range("I64").Value = 3106.4666
$ws.Range("J64").Value = 4297.086
$ws.Range("K64").Value = 3106.4666
$ws.Range("L64").Value = 4297.086
$ws.Range("M64").Value = -2858.4666
$ws.Range("N64").Value = -4793.086
$ws.Range("H67").Value = 3939.9
$ws.Range("I67").Value = 3106.4666
$ws.Range("J67").Value = 4297.086
$ws.Range("K67").Value = 3106.4666
$ws.Range("L67").Value = 4297.086
$ws.Range("M67").Value = -2248.4666
$ws.Range("N67").Value = -6013.086
$ws.Range("H88").Value = 2782
$ws.Range("I88").Value = 1934.3334
$ws.Range("J88").Value = 3099.875
$ws.Range("K88").Value = 1934.3334
$ws.Range("L88").Value = 3099.875
$ws.Range("M88").Value = -1528.3334
$ws.Range("N88").Value = -3911.875
$ws.Range("H91").Value = 2782
$ws.Range("I91").Value = 1934.3334
$ws.Range("J91").Value = 3099.875
$ws.Range("K91").Value = 1934.3334
$ws.Range("L91").Value = 3099.875
$ws.Range("M91").Value = -530.3334
$ws.Range("N91").Value = -5907.875
$ws.Range("H96").Value = 543.4286
$ws.Range("I96").Value = 507.33334
$ws.Range("J96").Value = 760
$ws.Range("K96").Value = 1522.00002
$ws.Range("L96").Value = 2280
$ws.Range("M96").Value = -149.0000199999999
$ws.Range("N96").Value = -5026
$ws.Range("H129").Value = 1000.8919
$ws.Range("J129").Value = 1184
$ws.Range("L129").Value = 3552
$ws.Range("N129").Value = -13552
$ws.Range("H132").Value = 1614.3671
$ws.Range("I132").Value = 1304.4058
$ws.Range("J132").Value = 3753.1
$ws.Range("K132").Value = 3913.2174
$ws.Range("L132").Value = 11259.3
$ws.Range("M132").Value = -1383.2174
$ws.Range("N132").Value = -16319.3
$ws.Range("H138").Value = 2820929.2
$ws.Range("I138").Value = 8001924
$ws.Range("J138").Value = 5171.087
$ws.Range("K138").Value = 24005772
$ws.Range("L138").Value = 15513.261
$ws.Range("M138").Value = -24000632
$ws.Range("N138").Value = -25793.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9417.615
$ws.Range("I32").Value = 9665.047
$ws.Range("J32").Value = 5161.8
$ws.Range("K32").Value = 9665.047
$ws.Range("L32").Value = 5161.8
$ws.Range("M32").Value = -9378.047
$ws.Range("N32").Value = -5735.8
$ws.Range("H61").Value = 1460.9546
$ws.Range("I61").Value = 1130.3158
$ws.Range("K61").Value = 1130.3158
$ws.Range("M61").Value = -918.3158000000001
$ws.Range("H63").Value = 3802.7856
$ws.Range("I63").Value = 3326.077
$ws.Range("K63").Value = 3326.077
$ws.Range("M63").Value = -2640.077
$ws.Range("H66").Value = 3802.7856
$ws.Range("I66").Value = 3326.077
$ws.Range("K66").Value = 16630.385
$ws.Range("M66").Value = -13198.385
$ws.Range("H92").Value = 158887.5
$ws.Range("J92").Value = 158887.5
$ws.Range("L92").Value = 158887.5
$ws.Range("N92").Value = -163879.5
$ws.Range("H97").Value = 802
$ws.Range("I97").Value = 802
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 802
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -306
$ws.Range("N97").ClearContents()
$ws.Range("H112").Value = 22971.5
$ws.Range("J112").Value = 22971.5
$ws.Range("L112").Value = 22971.5
$ws.Range("N112").Value = -25925.5
$ws.Range("H132").Value = 646311.7
$ws.Range("I132").Value = 769833.4
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 2309500.2
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -2306970.2
$ws.Range("N132").Value = -17057
$ws.Range("H136").Value = 1460.9546
$ws.Range("I136").Value = 1130.3158
$ws.Range("K136").Value = 3390.9474
$ws.Range("M136").Value = -840.9474

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1148.7931
$ws.Range("I94").Value = 950
$ws.Range("J94").Value = 1773.5714
$ws.Range("K94").Value = 950
$ws.Range("L94").Value = 1773.5714
$ws.Range("M94").Value = -499
$ws.Range("N94").Value = -2675.5714
$ws.Range("H105").Value = 2451.8823
$ws.Range("I105").Value = 2455.125
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 2455.125
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -708.125
$ws.Range("N105").Value = -5894
$ws.Range("H132").Value = 74164.914
$ws.Range("J132").Value = 74164.914
$ws.Range("L132").Value = 74164.914
$ws.Range("N132").Value = -84284.914

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 87.5
$ws.Range("I7").Value = 99.333336
$ws.Range("J7").Value = 52
$ws.Range("K7").Value = 99.333336
$ws.Range("L7").Value = 52
$ws.Range("M7").Value = 13.666664
$ws.Range("N7").Value = -278
$ws.Range("H31").Value = 12989288
$ws.Range("J31").Value = 3320.2
$ws.Range("L31").Value = 3320.2
$ws.Range("N31").Value = -3910.2
$ws.Range("H34").Value = 12989288
$ws.Range("J34").Value = 3320.2
$ws.Range("L34").Value = 3320.2
$ws.Range("N34").Value = -3724.2
$ws.Range("H134").Value = 509444.6
$ws.Range("I134").Value = 610930.7
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 1832792.1
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -1830257.1
$ws.Range("N134").Value = -11112
$ws.Range("H141").Value = 19624.5
$ws.Range("J141").Value = 19624.5
$ws.Range("L141").Value = 19624.5
$ws.Range("N141").Value = -29984.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 91.5
$ws.Range("I6").Value = 91.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 274.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -161.5
$ws.Range("N6").ClearContents()
$ws.Range("H99").Value = 5082.1665
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -754
$ws.Range("H113").Value = 844.8461
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 844.8461
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2534.5383
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6874.5383
$ws.Range("H120").Value = 10842.857
$ws.Range("I120").Value = 8750
$ws.Range("K120").Value = 26250
$ws.Range("M120").Value = -21412
$ws.Range("H131").Value = 850.89
$ws.Range("J131").Value = 884.43475
$ws.Range("L131").Value = 2653.30425
$ws.Range("N131").Value = -12733.30425
$ws.Range("H137").Value = 16669100
$ws.Range("I137").Value = 992.2222
$ws.Range("J137").Value = 30306644
$ws.Range("K137").Value = 2976.6666
$ws.Range("L137").Value = 90919932
$ws.Range("M137").Value = 2123.3334
$ws.Range("N137").Value = -90930132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10009
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10009
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10009
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -10549
$ws.Range("H73").Value = 10009
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10009
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10009
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -11881
$ws.Range("H80").Value = 3092.6924
$ws.Range("I80").Value = 2962.5
$ws.Range("J80").Value = 3301
$ws.Range("K80").Value = 2962.5
$ws.Range("L80").Value = 3301
$ws.Range("M80").Value = -1964.5
$ws.Range("N80").Value = -5297
$ws.Range("H83").Value = 3092.6924
$ws.Range("I83").Value = 2962.5
$ws.Range("J83").Value = 3301
$ws.Range("K83").Value = 14812.5
$ws.Range("L83").Value = 16505
$ws.Range("M83").Value = -9820.5
$ws.Range("N83").Value = -26489
$ws.Range("H102").Value = 2276.7036
$ws.Range("I102").Value = 2284.85
$ws.Range("J102").Value = 2253.4285
$ws.Range("K102").Value = 2284.85
$ws.Range("L102").Value = 2253.4285
$ws.Range("M102").Value = -662.8499999999999
$ws.Range("N102").Value = -5497.4285
$ws.Range("H113").Value = 1090.5294
$ws.Range("I113").Value = 963.9167
$ws.Range("J113").Value = 1394.4
$ws.Range("K113").Value = 963.9167
$ws.Range("L113").Value = 1394.4
$ws.Range("M113").Value = 1206.0833
$ws.Range("N113").Value = -5734.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 1581.5
$ws.Range("J9").Value = 1741.8
$ws.Range("L9").Value = 1741.8
$ws.Range("N9").Value = -2189.8
$ws.Range("H40").Value = 6913.3335
$ws.Range("I40").Value = 7696
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 7696
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -7560
$ws.Range("N40").Value = -3272
$ws.Range("H68").Value = 1228.7858
$ws.Range("I68").Value = 930.6957
$ws.Range("J68").Value = 2600
$ws.Range("K68").Value = 930.6957
$ws.Range("L68").Value = 2600
$ws.Range("M68").Value = -181.6957
$ws.Range("N68").Value = -4098
$ws.Range("H71").Value = 1228.7858
$ws.Range("I71").Value = 930.6957
$ws.Range("J71").Value = 2600
$ws.Range("K71").Value = 4653.4785
$ws.Range("L71").Value = 13000
$ws.Range("M71").Value = -909.4785000000002
$ws.Range("N71").Value = -20488
$ws.Range("H96").Value = 30197
$ws.Range("J96").Value = 30197
$ws.Range("L96").Value = 30197
$ws.Range("N96").Value = -35689

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H96").Value = 1200
$ws.Range("J96").Value = 1200
$ws.Range("L96").Value = 1200
$ws.Range("N96").Value = -3946
$ws.Range("H136").Value = 1692.8043
$ws.Range("I136").Value = 1744.3636
$ws.Range("J136").Value = 1561.9231
$ws.Range("K136").Value = 5233.0908
$ws.Range("L136").Value = 4685.7693
$ws.Range("M136").Value = -2683.0908
$ws.Range("N136").Value = -9785.7693
